$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "time_taken", matching the header style used by B1:E1
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Data cells F2:F58: per-row "time_taken" timestamp values
$timeTaken = @(
    "2021-10-05 10:52:50.275921",
    "2021-10-05 10:52:50.275934",
    "2021-10-05 10:52:50.275938",
    "2021-10-05 10:52:50.275941",
    "2021-10-05 10:52:50.275944",
    "2021-10-05 10:52:50.275947",
    "2021-10-05 10:52:50.275950",
    "2021-10-05 10:52:50.275953",
    "2021-10-05 10:52:50.275957",
    "2021-10-05 10:52:50.275960",
    "2021-10-05 10:52:50.275963",
    "2021-10-05 10:52:50.275966",
    "2021-10-05 10:52:50.275969",
    "2021-10-05 10:52:50.275972",
    "2021-10-05 10:52:50.275975",
    "2021-10-05 10:52:50.275978",
    "2021-10-05 10:52:50.275981",
    "2021-10-05 10:52:50.275984",
    "2021-10-05 10:52:50.275987",
    "2021-10-05 10:52:50.275990",
    "2021-10-05 10:52:50.275993",
    "2021-10-05 10:52:50.275996",
    "2021-10-05 10:52:50.275999",
    "2021-10-05 10:52:50.276002",
    "2021-10-05 10:52:50.276005",
    "2021-10-05 10:52:50.276008",
    "2021-10-05 10:52:50.276011",
    "2021-10-05 10:52:50.276014",
    "2021-10-05 10:52:50.276017",
    "2021-10-05 10:52:50.276020",
    "2021-10-05 10:52:50.276023",
    "2021-10-05 10:52:50.276026",
    "2021-10-05 10:52:50.276030",
    "2021-10-05 10:52:50.276033",
    "2021-10-05 10:52:50.276036",
    "2021-10-05 10:52:50.276039",
    "2021-10-05 10:52:50.276042",
    "2021-10-05 10:52:50.276044",
    "2021-10-05 10:52:50.276047",
    "2021-10-05 10:52:50.276050",
    "2021-10-05 10:52:50.276054",
    "2021-10-05 10:52:50.276057",
    "2021-10-05 10:52:50.276060",
    "2021-10-05 10:52:50.276063",
    "2021-10-05 10:52:50.276065",
    "2021-10-05 10:52:50.276068",
    "2021-10-05 10:52:50.276071",
    "2021-10-05 10:52:50.276074",
    "2021-10-05 10:52:50.276077",
    "2021-10-05 10:52:50.276080",
    "2021-10-05 10:52:50.276083",
    "2021-10-05 10:52:50.276086",
    "2021-10-05 10:52:50.276089",
    "2021-10-05 10:52:50.276092",
    "2021-10-05 10:52:50.276096",
    "2021-10-05 10:52:50.276098",
    "2021-10-05 10:52:50.276101"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

